$wb = $excel.ActiveWorkbook

# --- Rename the three activity sheets (3/4/5 -> 6/7/8) ---
$wb.Worksheets.Item(1).Name = "Activité 6"
$wb.Worksheets.Item(2).Name = "Activité 7"
$wb.Worksheets.Item(3).Name = "Activité 8"

# --- Update header/footer font style on every sheet: "Regular" -> "Normal" ---
foreach ($ws in $wb.Worksheets) {
    $ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
    $ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
}

# --- Move the active tab from the 3rd sheet to the 1st sheet ---
$wb.Worksheets.Item(1).Activate()

Write-Host "Done"
